# Weekly automated refresh of the cryptocurrency price/volume snapshot
# (GitHub Actions job). Updates the Price (D) and Volume(1h) (E) columns
# for every listed coin; rows 43/44 also swap ranking position (Cronos
# now sits above HuobiToken), so their Coin name + Link columns change too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on price cells whose new value looks numeric
# (prevents Excel from auto-converting them to actual numbers, which
# would both change their type and mangle trailing zeros / precision).
$textCells = @("D4", "D5", "D6", "D8", "D10", "D12", "D14", "D15", "D16", "D19", "D20", "D22", "D24", "D25", "D26", "D27", "D29", "D32", "D33", "D34", "D35", "D37", "D39", "D40", "D41", "D43", "D44", "D46", "D47", "D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "37.365.07"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").Value = "2.054.73"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "231.05"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("D6").Value = "0.622"
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "57.25"
$ws.Range("E8").Value = "  -3.30%  "
$ws.Range("E9").Value = "  -2.03%  "
$ws.Range("D10").Value = "0.0771"
$ws.Range("E10").Value = "  -2.49%  "
$ws.Range("E11").Value = "  +1.41%  "
$ws.Range("D12").Value = "14.78"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "2.354.89"
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("D14").Value = "20.66"
$ws.Range("E14").Value = "  -2.59%  "
$ws.Range("D15").Value = "0.759"
$ws.Range("E15").Value = "  -1.90%  "
$ws.Range("D16").Value = "5.26"
$ws.Range("E16").Value = "  -1.65%  "
$ws.Range("D17").Value = "2.063.68"
$ws.Range("E17").Value = "  -3.59%  "
$ws.Range("D18").Value = "37.359.71"
$ws.Range("E18").Value = "  -0.96%  "
$ws.Range("D19").Value = "6.03"
$ws.Range("E19").Value = "  -1.90%  "
$ws.Range("D20").Value = "69.89"
$ws.Range("E20").Value = "  -2.37%  "
$ws.Range("D21").Value = "0.0₃0824"
$ws.Range("E21").Value = "  -3.25%  "
$ws.Range("D22").Value = "227.22"
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Value = "2.39"
$ws.Range("E24").Value = "  +1.03%  "
$ws.Range("D25").Value = "2.34"
$ws.Range("E25").Value = "  -3.45%  "
$ws.Range("D26").Value = "9.59"
$ws.Range("E26").Value = "  +4.53%  "
$ws.Range("D27").Value = "169.36"
$ws.Range("E27").Value = "  -0.63%  "
$ws.Range("E28").Value = "  -1.92%  "
$ws.Range("D29").Value = "19.18"
$ws.Range("E29").Value = "  -1.55%  "
$ws.Range("E30").Value = "  -4.83%  "
$ws.Range("E31").Value = "  +0.54%  "
$ws.Range("D32").Value = "4.53"
$ws.Range("E32").Value = "  -3.93%  "
$ws.Range("D33").Value = "0.0625"
$ws.Range("E33").Value = "  -1.19%  "
$ws.Range("D34").Value = "4.58"
$ws.Range("E34").Value = "  -3.14%  "
$ws.Range("D35").Value = "2.50"
$ws.Range("E35").Value = "  +0.55%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").Value = "3.29"
$ws.Range("E37").Value = "  -3.36%  "
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("D39").Value = "5.27"
$ws.Range("E39").Value = "  -2.46%  "
$ws.Range("D40").Value = "0.0226"
$ws.Range("E40").Value = "  +4.49%  "
$ws.Range("D41").Value = "98.76"
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("D42").Value = "1.488.70"
$ws.Range("E42").Value = "  +3.17%  "
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").Value = "0.0954"
$ws.Range("E43").Value = "  -3.06%  "
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").Value = "2.90"
$ws.Range("E44").Value = "  +0.61%  "
$ws.Range("E45").Value = "  +2.69%  "
$ws.Range("D46").Value = "16.74"
$ws.Range("E46").Value = "  +0.75%  "
$ws.Range("D47").Value = "4.00"
$ws.Range("E47").Value = "  -6.17%  "
$ws.Range("E48").Value = "  -2.91%  "
$ws.Range("D49").Value = "7.26"
$ws.Range("E49").Value = "  -1.83%  "
$ws.Range("E50").Value = "  -2.48%  "
$ws.Range("D51").Value = "2.240.24"
$ws.Range("E51").Value = "  -1.23%  "

# Restore default (General) cell style now that the text values are
# committed, so the cells end up styled exactly like before the edit.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
